# Update "想去人数" (F column) values on the sheets that contain data rows:
# "展览" (Worksheet 1) and "全部类型" (Worksheet 4)

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1315
    $ws.Range("F3").Value = 1734
    $ws.Range("F5").Value = 6264
    $ws.Range("F6").Value = 99
}
